$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "European Commission's Humanitarian Aid and Civil Protection Department"
$ws.Range("B2").Value = 6063308

$ws.Range("A3").Value = "Central Emergency Response Fund"
$ws.Range("B3").Value = 3049997

$ws.Range("A4").Value = "UN COVID-19 Response and Recovery Fund"
$ws.Range("B4").Value = 893000

$ws.Range("A5").Value = "Japan, Government of"
$ws.Range("B5").Value = 2121418

$ws.Range("A6").Value = "Germany, Government of"
$ws.Range("B6").Value = 19514
